{"js": "// Insert a new \"Recommendation Type\" / \"Compressor\" row at the top of the\n// \"Summary of Estimated Savings and Implementation Costs\" table, immediately\n// before the existing \"Annual Cost Savings\" row.\n//\n// Locate that table via the existing \"Annual Cost Savings\" row (more robust\n// than assuming a fixed table index).\nconst results = context.document.body.search(\"Annual Cost Savings\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nconst table = results.items[0].parentTable;\ntable.addRows(Word.InsertLocation.start, 1, [[\"Recommendation Type\", \"Compressor\"]]);\nawait context.sync();\n", "ps1": "# Insert a new \"Recommendation Type\" / \"Compressor\" row at the top of the\n# \"Summary of Estimated Savings and Implementation Costs\" table, immediately\n# before the existing \"Annual Cost Savings\" row.\n#\n# Locate that table via the existing \"Annual Cost Savings\" row (more robust\n# than assuming a fixed table index).\n$d = $word.ActiveDocument\n$range = $d.Content\n$range.Find.ClearFormatting()\n$found = $range.Find.Execute(\"Annual Cost Savings\")\n$table = $range.Tables.Item(1)\n\n$newRow = $table.Rows.Add($table.Rows.Item(1))\n$newRow.Cells.Item(1).Range.Text = \"Recommendation Type\"\n$newRow.Cells.Item(2).Range.Text = \"Compressor\"\n"}
